# APfollowup.xlsx update — add three new follow-up items to Sheet1,
# turn on word-wrap for column B, and grow row 7 to fit its wrapped text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 7 (item 6) gets its description filled in, and two brand
# new rows (items 7 and 8) are appended below it. Values are written in
# this particular order so new shared-string entries land at indices
# 12, 13, 14 in the same order as the authored workbook.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "准备思维导图，想想处理逻辑"

$ws.Range("B7").Value = "先从杭州做起：浙江省高考政策，高考学生数据，高考学生人数，坐落于杭州的大学"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "先做微信版本，APP先放放。"

# Enable wrap text on the whole of column B (adds a new wrapText cell
# style and applies it to every populated cell in the column).
$ws.Columns("B").WrapText = $true

# Row 7 now holds a long wrapped line, so grow it to fit.
$ws.Rows(7).RowHeight = 28.5

# Selection moves on to the next empty row, ready for the following entry.
$ws.Range("B10").Select() | Out-Null
